$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.547.86"
$ws.Range("E2").Value = "  -1.30%  "
$ws.Range("D3").Value = "2.996.81"
$ws.Range("E3").Value = "  -0.40%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.24"
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.74"
$ws.Range("E6").Value = "  -1.65%  "
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.521"
$ws.Range("E8").Value = "  +0.44%  "
$ws.Range("D9").Value = "2.994.19"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("E10").Value = "  -1.04%  "
$ws.Range("E11").Value = "  +6.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.457"
$ws.Range("E12").Value = "  +3.68%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000229"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.44"
$ws.Range("E14").Value = "  -0.78%  "
$ws.Range("E15").Value = "  +2.19%  "
$ws.Range("D16").Value = "3.492.90"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.98"
$ws.Range("E17").Value = "  -0.74%  "
$ws.Range("D18").Value = "61.482.94"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "2.999.81"
$ws.Range("E19").Value = "  -0.40%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "451.91"
$ws.Range("E20").Value = "  -1.53%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.04"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.686"
$ws.Range("E22").Value = "  +0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.34"
$ws.Range("E23").Value = "  +0.30%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "81.64"
$ws.Range("E24").Value = "  +2.11%  "
$ws.Range("B25").Value = "RenderToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "10.80"
$ws.Range("E25").Value = "  +6.63%  "
$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("E26").Value = "  -2.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.96"
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("E29").Value = "  +2.90%  "
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.21"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("E32").Value = "  -1.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.34"
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  +3.21%  "
$ws.Range("D35").Value = "0.0₃0828"
$ws.Range("E35").Value = "  +5.17%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.79"
$ws.Range("E37").Value = "  +0.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "9.22"
$ws.Range("E38").Value = "  +2.51%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "50.44"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.07"
$ws.Range("E40").Value = "  -2.11%  "
$ws.Range("E41").Value = "  +11.00%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  -0.63%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "398.00"
$ws.Range("E43").Value = "  -4.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.67"
$ws.Range("E44").Value = "  +4.34%  "
$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.271"
$ws.Range("E45").Value = "  -1.43%  "
$ws.Range("B46").Value = "VeChain"
$ws.Range("C46").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0354"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("D47").Value = "2.694.47"
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "131.18"
$ws.Range("E48").Value = "  +1.97%  "
$ws.Range("E49").Value = "  +0.12%  "
$ws.Range("E50").Value = "  -0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.15"
$ws.Range("E51").Value = "  +1.71%  "
